$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '26.566.05'
$ws.Cells.Item(2, 5).Value = '  -1.78%  '

$ws.Cells.Item(3, 4).Value = '1.795.38'
$ws.Cells.Item(3, 5).Value = '  -0.20%  '

$ws.Cells.Item(4, 4).Value = "'1.000"
$ws.Cells.Item(4, 5).Value = '  -0.09%  '

$ws.Cells.Item(5, 2).Value = 'BNB'
$ws.Cells.Item(5, 3).Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Cells.Item(5, 4).Value = "'309.19"
$ws.Cells.Item(5, 5).Value = '  +0.57%  '

$ws.Cells.Item(6, 2).Value = 'USDC'
$ws.Cells.Item(6, 3).Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Cells.Item(6, 4).Value = "'1.000"
$ws.Cells.Item(6, 5).Value = '  -0.03%  '

$ws.Cells.Item(7, 4).Value = "'0.4284"
$ws.Cells.Item(7, 5).Value = '  +1.97%  '

$ws.Cells.Item(8, 4).Value = "'0.3627"
$ws.Cells.Item(8, 5).Value = '  +1.31%  '

$ws.Cells.Item(9, 4).Value = "'0.07176"
$ws.Cells.Item(9, 5).Value = '  +0.80%  '

$ws.Cells.Item(10, 4).Value = "'0.8571"
$ws.Cells.Item(10, 5).Value = '  +1.38%  '

$ws.Cells.Item(11, 4).Value = "'20.63"
$ws.Cells.Item(11, 5).Value = '  +2.25%  '

$ws.Cells.Item(12, 4).Value = '1.894.12'
$ws.Cells.Item(12, 5).Value = '  +9.72%  '

$ws.Cells.Item(13, 4).Value = "'6.531"
$ws.Cells.Item(13, 5).Value = '  +2.81%  '

$ws.Cells.Item(14, 4).Value = "'5.291"
$ws.Cells.Item(14, 5).Value = '  -0.13%  '

$ws.Cells.Item(15, 4).Value = "'0.06900"
$ws.Cells.Item(15, 5).Value = '  +1.80%  '

$ws.Cells.Item(16, 4).Value = "'1.000"
$ws.Cells.Item(16, 5).Value = '  -0.44%  '

$ws.Cells.Item(17, 5).Value = '  -0.64%  '

$ws.Cells.Item(18, 4).Value = "'0.000008776"
$ws.Cells.Item(18, 5).Value = '  +0.86%  '

$ws.Cells.Item(19, 4).Value = "'1.002"
$ws.Cells.Item(19, 5).Value = '  +0.18%  '

$ws.Cells.Item(20, 4).Value = "'15.04"
$ws.Cells.Item(20, 5).Value = '  -0.03%  '

$ws.Cells.Item(21, 4).Value = '26.574.97'
$ws.Cells.Item(21, 5).Value = '  -1.77%  '

$ws.Cells.Item(22, 4).Value = "'5.139"

$ws.Cells.Item(23, 5).Value = '  +0.28%  '

$ws.Cells.Item(24, 4).Value = '2.114.60'
$ws.Cells.Item(24, 5).Value = '  +6.53%  '

$ws.Cells.Item(25, 4).Value = "'151.93"
$ws.Cells.Item(25, 5).Value = '  -0.60%  '

$ws.Cells.Item(26, 4).Value = "'1.824"
$ws.Cells.Item(26, 5).Value = '  -5.58%  '

$ws.Cells.Item(27, 4).Value = "'18.19"
$ws.Cells.Item(27, 5).Value = '  +0.20%  '

$ws.Cells.Item(28, 4).Value = "'5.164"
$ws.Cells.Item(28, 5).Value = '  +2.86%  '

$ws.Cells.Item(29, 4).Value = "'1.887"
$ws.Cells.Item(29, 5).Value = '  +14.68%  '

$ws.Cells.Item(30, 5).Value = '  +1.73%  '

$ws.Cells.Item(31, 4).Value = "'0.08903"
$ws.Cells.Item(31, 5).Value = '  -1.12%  '

$ws.Cells.Item(32, 4).Value = "'0.7469"
$ws.Cells.Item(32, 5).Value = '  +3.32%  '

$ws.Cells.Item(33, 4).Value = "'1.149"
$ws.Cells.Item(33, 5).Value = '  +5.71%  '

$ws.Cells.Item(34, 4).Value = "'4.367"
$ws.Cells.Item(34, 5).Value = '  +1.43%  '

$ws.Cells.Item(35, 2).Value = 'Frax'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Cells.Item(35, 4).Value = "'1.002"
$ws.Cells.Item(35, 5).Value = '  +0.19%  '

$ws.Cells.Item(36, 2).Value = 'HuobiToken'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Cells.Item(36, 4).Value = "'2.740"
$ws.Cells.Item(36, 5).Value = '  -4.18%  '

$ws.Cells.Item(37, 4).Value = "'1.110"
$ws.Cells.Item(37, 5).Value = '  +2.88%  '

$ws.Cells.Item(38, 5).Value = '  +0.68%  '

$ws.Cells.Item(39, 4).Value = "'0.01906"
$ws.Cells.Item(39, 5).Value = '  +0.03%  '

$ws.Cells.Item(40, 4).Value = "'0.4993"
$ws.Cells.Item(40, 5).Value = '  +0.71%  '

$ws.Cells.Item(41, 4).Value = "'0.1624"
$ws.Cells.Item(41, 5).Value = '  -0.15%  '

$ws.Cells.Item(42, 4).Value = "'2.615"
$ws.Cells.Item(42, 5).Value = '  +0.87%  '

$ws.Cells.Item(43, 4).Value = "'6.446"
$ws.Cells.Item(43, 5).Value = '  +8.91%  '

$ws.Cells.Item(44, 4).Value = "'8.243"
$ws.Cells.Item(44, 5).Value = '  +2.48%  '

$ws.Cells.Item(45, 4).Value = "'105.89"
$ws.Cells.Item(45, 5).Value = '  +0.94%  '

$ws.Cells.Item(46, 4).Value = "'10.33"
$ws.Cells.Item(46, 5).Value = '  +1.61%  '

$ws.Cells.Item(47, 4).Value = "'1.001"
$ws.Cells.Item(47, 5).Value = '  +0.05%  '

$ws.Cells.Item(48, 4).Value = "'1.643"
$ws.Cells.Item(48, 5).Value = '  +2.48%  '

$ws.Cells.Item(49, 2).Value = 'Decentraland'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Cells.Item(49, 4).Value = "'0.4509"
$ws.Cells.Item(49, 5).Value = '  -0.49%  '

$ws.Cells.Item(50, 2).Value = 'Cronos'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Cells.Item(50, 4).Value = "'0.06211"
$ws.Cells.Item(50, 5).Value = '  -1.38%  '

$ws.Cells.Item(51, 4).Value = "'1.782"
$ws.Cells.Item(51, 5).Value = '  +4.36%  '
